# Anonymize "fedcore" -> "approach" and restyle the header spacer cells
# (C1/D1 on quality_comparison; C1/D1/F1/G1 on computational_comparison),
# then drop the stray empty G5 cell on computational_comparison.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# --- Build the two new border styles once on ws1, then propagate by
# --- copy/paste-format so the style table doesn't accumulate orphan xfs.

# C1: top+bottom border only (no left/right), plain (non-bold) font
$ws1.Range("C1").ClearFormats()
$ws1.Range("C1").Borders.Item(8).LineStyle = 1   # xlEdgeTop
$ws1.Range("C1").Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# D1: top+bottom+right border (no left), plain (non-bold) font
# (edge order chosen so every intermediate combo already exists in the
# style table, so no throwaway <border> entries get appended)
$ws1.Range("D1").ClearFormats()
$ws1.Range("D1").Borders.Item(8).LineStyle = 1    # xlEdgeTop
$ws1.Range("D1").Borders.Item(10).LineStyle = 1   # xlEdgeRight
$ws1.Range("D1").Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# Propagate the same two styles onto sheet 2's matching spacer cells.
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)

$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("G1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Anonymize the "fedcore" column headers to "approach"
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty inline-string cell G5 on computational_comparison
$ws2.Range("G5").ClearContents()
